$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that originally sat
#    right under the "Play Beellionaires Dream Drop Free - Review"
#    Heading1 paragraph.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2. At the end of the document, turn the single DALLE image-prompt
#    paragraph into two paragraphs:
#       - a bold "Play Beellionaires Dream Drop Free - Review" line
#       - the (still italic) meta-description sentence
#    We do this with one InsertXML call that rewrites the whole last
#    paragraph's contents into the two target paragraphs, then clean up
#    the left-over empty paragraph mark that InsertXML leaves behind.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastRange = $lastPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Beellionaires Dream Drop Free - Review</w:t></w:r></w:p>' +
  '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Beellionaires Dream Drop - a Relax Gaming slot game featuring potential winnings up to 10,000 times your bet. Play the game for free.</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$lastRange.InsertXML($xml)

# InsertXML leaves an extra empty paragraph mark behind (because it
# replaces the range's *content* but the paragraph's own end-of-paragraph
# mark survives as a stray empty paragraph). Remove it.
$newCount = $d.Paragraphs.Count
$trailing = $d.Paragraphs($newCount)
$trailingRange = $trailing.Range
$prevEnd = $d.Paragraphs($newCount - 1).Range.End
$cleanupRange = $d.Range($prevEnd - 1, $trailingRange.End)
$cleanupRange.Delete()
